$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.650.00'
$ws.Range("E2").Value = '  -5.79%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.901.19'
$ws.Range("E3").Value = '  -3.62%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.24'
$ws.Range("E5").Value = '  -2.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '123.33'
$ws.Range("E6").Value = '  -4.09%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.899.55'
$ws.Range("E8").Value = '  -3.65%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.125'
$ws.Range("E10").Value = '  -7.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.74'
$ws.Range("E11").Value = '  -9.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.439'
$ws.Range("E12").Value = '  +2.16%  '

$ws.Range("E13").Value = '  -4.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.38'
$ws.Range("E14").Value = '  -1.13%  '

$ws.Range("E15").Value = '  +1.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.372.80'
$ws.Range("E16").Value = '  -3.89%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.890.95'
$ws.Range("E17").Value = '  -4.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.58'
$ws.Range("E18").Value = '  +5.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '57.602.09'
$ws.Range("E19").Value = '  -6.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '407.48'
$ws.Range("E20").Value = '  -6.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.92'
$ws.Range("E21").Value = '  -1.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.671'
$ws.Range("E22").Value = '  +1.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.84'
$ws.Range("E23").Value = '  -4.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.91'
$ws.Range("E24").Value = '  +3.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '77.13'
$ws.Range("E25").Value = '  -2.32%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.16%  '

$ws.Range("E28").Value = '  -1.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.94'
$ws.Range("E29").Value = '  +3.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.24'
$ws.Range("E30").Value = '  +0.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.05'
$ws.Range("E31").Value = '  -2.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '24.70'
$ws.Range("E32").Value = '  -3.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0985'
$ws.Range("E33").Value = '  +5.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.912'
$ws.Range("E34").Value = '  -4.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.39'
$ws.Range("E35").Value = '  -2.49%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.02'
$ws.Range("E36").Value = '  -10.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '48.22'
$ws.Range("E37").Value = '  -3.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.47'
$ws.Range("E38").Value = '  +9.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0625'
$ws.Range("E39").Value = '  -7.55%  '

$ws.Range("E40").Value = '  -5.33%  '

$ws.Range("E41").Value = '  -1.44%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.44'
$ws.Range("E42").Value = '  +0.50%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.616.29'
$ws.Range("E43").Value = '  -1.32%  '

$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '359.71'
$ws.Range("E44").Value = '  -3.39%  '

$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '120.13'
$ws.Range("E46").Value = '  +0.57%  '

$ws.Range("E47").Value = '  -2.73%  '

$ws.Range("E48").Value = '  +0.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.94'
$ws.Range("E49").Value = '  -1.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.94'
$ws.Range("E50").Value = '  -2.53%  '

$ws.Range("E51").Value = '  -3.83%  '
